$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Row 3: playerID value changed (stays text via leading apostrophe), operatorID (betvita->10909001)
# loses its explicit style, and the player/platform token gets a fresh value.
$ws.Range("A3").Value = "'1000041"

$ws.Range("C3").Value = "'10909001"
$ws.Range("C3").Style = "Normal"

$ws.Range("D3").Value = "c0987b11-c5d6-4746-a763-ac07fde2e4b8-1731498207823"

# Row 7: new authToken value (A7) and a new response-validation token (B7),
# with B7 losing its previous explicit (blank) styling.
$ws.Range("A7").Value = "d628b6d5-cfe0-41d1-89c4-d4964e76995c"

$ws.Range("B7").Value = "2712dd32-8386-4610-b361-ebb1329f0aa5"
$ws.Range("B7").Style = "Normal"

# Update the active selection on the sheet to D3, matching the saved view state.
$ws.Activate()
$ws.Range("D3").Select()
